# Commit message: "replace buffer.name with buffer.id"
#
# The "buffer" worksheet used to have a first column "name" that stored a
# concatenation like "<item> @ <location>". That column is no longer
# needed (the sheet is keyed by type/location/item instead), so the
# entire column A is deleted from the "buffer" sheet, shifting the
# remaining columns (type, location, item) one position to the left.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("buffer")

# Delete the whole "name" column (column A) from the buffer sheet.
$ws.Columns.Item(1).Delete()
